$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 828
$ws.Range("B7").Value = 0.04252751383954506

$ws.Range("A14").Value = 744
$ws.Range("B14").Value = 0.1682054113418063

$ws.Range("A21").Value = 616
$ws.Range("B21").Value = 0.229812642500469
